# Daily attendance processing - 2025-10-07 12:35:49
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder "Recorded By" email lists (same set of recipients, new order) ---
$ws.Range("G3").Value  = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("G25").Value = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

$ws.Range("G12").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G34").Value = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

$ws.Range("G41").Value = "marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

# --- Class / group statistics numbers updated by the daily processing run ---
$ws.Range("L7").Value = 3
$ws.Range("L8").Value = 36

$ws.Range("P15").Value = 2
$ws.Range("Q15").Value = 18

# --- Row 11 (Biochemistry Lab/CBL session 1) got recorded/processed:
#     status flips from "Pending" to "Not Recorded", and the row's
#     highlight color switches from the "Pending" (yellow) style to the
#     "Not Recorded" (pink) style used elsewhere in the sheet (e.g. row 2).
$ws.Range("A2:I2").Copy()
$ws.Range("A11:I11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I11").Value = "Not Recorded"
